# Spain Segunda workbook: corrects row-ordering within several same-date match
# blocks (the "id" column and all stats, B:AD, were attached to the wrong
# A-index row) - row content (columns B through AD) is permuted between rows
# while column A (the sequential match index) stays put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot the "before" B:AD content of every affected row first, so that the
#    swaps/cycles below always read original values, never a value already
#    overwritten earlier in the script.
$buf14 = $ws.Range("B14:AD14").Value()
$buf15 = $ws.Range("B15:AD15").Value()
$buf36 = $ws.Range("B36:AD36").Value()
$buf37 = $ws.Range("B37:AD37").Value()
$buf45 = $ws.Range("B45:AD45").Value()
$buf46 = $ws.Range("B46:AD46").Value()
$buf49 = $ws.Range("B49:AD49").Value()
$buf50 = $ws.Range("B50:AD50").Value()
$buf96 = $ws.Range("B96:AD96").Value()
$buf97 = $ws.Range("B97:AD97").Value()
$buf169 = $ws.Range("B169:AD169").Value()
$buf170 = $ws.Range("B170:AD170").Value()
$buf185 = $ws.Range("B185:AD185").Value()
$buf186 = $ws.Range("B186:AD186").Value()
$buf190 = $ws.Range("B190:AD190").Value()
$buf191 = $ws.Range("B191:AD191").Value()
$buf213 = $ws.Range("B213:AD213").Value()
$buf214 = $ws.Range("B214:AD214").Value()
$buf220 = $ws.Range("B220:AD220").Value()
$buf221 = $ws.Range("B221:AD221").Value()
$buf230 = $ws.Range("B230:AD230").Value()
$buf231 = $ws.Range("B231:AD231").Value()
$buf247 = $ws.Range("B247:AD247").Value()
$buf248 = $ws.Range("B248:AD248").Value()
$buf371 = $ws.Range("B371:AD371").Value()
$buf372 = $ws.Range("B372:AD372").Value()
$buf449 = $ws.Range("B449:AD449").Value()
$buf450 = $ws.Range("B450:AD450").Value()
$buf451 = $ws.Range("B451:AD451").Value()
$buf452 = $ws.Range("B452:AD452").Value()
$buf453 = $ws.Range("B453:AD453").Value()
$buf454 = $ws.Range("B454:AD454").Value()
$buf455 = $ws.Range("B455:AD455").Value()
$buf456 = $ws.Range("B456:AD456").Value()
$buf457 = $ws.Range("B457:AD457").Value()

# 2) Write the permuted content back: row R receives the snapshot that originally
#    belonged to its paired/rotated row.
$ws.Range("B14:AD14").Value = $buf15
$ws.Range("B15:AD15").Value = $buf14
$ws.Range("B36:AD36").Value = $buf37
$ws.Range("B37:AD37").Value = $buf36
$ws.Range("B45:AD45").Value = $buf46
$ws.Range("B46:AD46").Value = $buf45
$ws.Range("B49:AD49").Value = $buf50
$ws.Range("B50:AD50").Value = $buf49
$ws.Range("B96:AD96").Value = $buf97
$ws.Range("B97:AD97").Value = $buf96
$ws.Range("B169:AD169").Value = $buf170
$ws.Range("B170:AD170").Value = $buf169
$ws.Range("B185:AD185").Value = $buf186
$ws.Range("B186:AD186").Value = $buf185
$ws.Range("B190:AD190").Value = $buf191
$ws.Range("B191:AD191").Value = $buf190
$ws.Range("B213:AD213").Value = $buf214
$ws.Range("B214:AD214").Value = $buf213
$ws.Range("B220:AD220").Value = $buf221
$ws.Range("B221:AD221").Value = $buf220
$ws.Range("B230:AD230").Value = $buf231
$ws.Range("B231:AD231").Value = $buf230
$ws.Range("B247:AD247").Value = $buf248
$ws.Range("B248:AD248").Value = $buf247
$ws.Range("B371:AD371").Value = $buf372
$ws.Range("B372:AD372").Value = $buf371
$ws.Range("B449:AD449").Value = $buf453
$ws.Range("B450:AD450").Value = $buf449
$ws.Range("B451:AD451").Value = $buf450
$ws.Range("B452:AD452").Value = $buf451
$ws.Range("B453:AD453").Value = $buf452
$ws.Range("B454:AD454").Value = $buf455
$ws.Range("B455:AD455").Value = $buf456
$ws.Range("B456:AD456").Value = $buf457
$ws.Range("B457:AD457").Value = $buf454
